$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1. "Benedict Gaster" -> "Benedict " + proofErr(spellStart) + "Gaster" + proofErr(spellEnd)
# ---------------------------------------------------------------------
$rngName = $d.Content
$foundName = $rngName.Find.Execute("Benedict Gaster", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundName) {
    $xmlName = "<w:p $wns><w:r><w:t xml:space=`"preserve`">Benedict </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Gaster</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"
    $rngName.InsertXML($xmlName)
}

# ---------------------------------------------------------------------
# 2. Fill in the first blank meeting-log row (Date / Notes / Actions)
# ---------------------------------------------------------------------
$t = $d.Tables.Item(2)
$row = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $candidate = $t.Rows.Item($i)
    $c1text = $candidate.Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)
    $c2text = $candidate.Cells.Item(2).Range.Text.TrimEnd([char]13, [char]7)
    if ($c1text -eq "" -and $c2text -eq "") {
        $row = $candidate
        break
    }
}

# Cell 1: Date
$row.Cells.Item(1).Range.Text = "11/Nov/2021"

# Cell 3: Actions
$row.Cells.Item(3).Range.Text = "Add things discussed to literature review"

# Cell 2: Meeting notes (several paragraphs)
$cell2 = $row.Cells.Item(2)
$lines = @(
  "Talked about literature review:",
  "Each chapter should have an intro, what it contains",
  "Game engines",
  "VR",
  "Web tech, frontend and backend",
  "Why using Django rather than flask",
  "Oculus integrations vs XR",
  "Rest api Django has these capabilities "
)
$cell2.Range.Text = [string]::Join("`r", $lines)

# Fix paragraph 1 of cell 2: split "Talked about literature review:" into two runs
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Talked about literature review:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $xml1 = "<w:p $wns><w:r><w:t>Talked about literature review</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>"
    $rng1.InsertXML($xml1)
}

# Fix last paragraph of cell 2: mark "api" with spell-check proofErr
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Rest api Django has these capabilities ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $xml2 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">Rest </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> Django has these capabilities </w:t></w:r></w:p>"
    $rng2.InsertXML($xml2)
}
